$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab.
$ws.Name = "g vs. hardcoded C 6 checkers"

# 2. Append rows 38-47 with the same look (bold/bordered/centered column A,
#    text-looking numbers in A, plain numbers in B/C) as the existing rows.
$newRows = @(
    @(38, "56", 0,   33),
    @(39, "57", 100, 21),
    @(40, "58", 0,   3),
    @(41, "59", 100, 24),
    @(42, "60", 0,   1),
    @(43, "61", 0,   33),
    @(44, "62", 100, 28),
    @(45, "63", 0,   33),
    @(46, "64", 0,   33),
    @(47, "65", 0,   33)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $aText = $row[1]
    $bVal = $row[2]
    $cVal = $row[3]

    # Copy the formatting of the last existing data row (row 37) into column
    # A of the new row so the new label cell keeps the bold/border/center
    # style already used for every other row in that column.
    $ws.Cells.Item(37, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null

    # Leading apostrophe forces the numeric-looking label to be stored as
    # text, matching the existing column A cells.
    $ws.Cells.Item($r, 1).Value = "'" + $aText

    $ws.Cells.Item($r, 2).Value = $bVal
    $ws.Cells.Item($r, 3).Value = $cVal
}

$excel.CutCopyMode = $false
